# Commit: "renamed storage::manager to storage::api and *Manager.java to *Storage.java"
#
# Applies to the single "StorageComponent" diagram slide:
#   - storage::manager  -> storage::api   (rename the storage API facade label)
#   - EvaluationsManager -> EvaluationsStorage
#   - AccountsManager    -> AccountsStorage
#   - CoursesManager      -> CoursesStorage
#   - nudge the connector under the storage::api box to line up with its
#     (now shorter) label

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "storage::manager" -> "storage::api" -------------------------------
# This text lives in the "Rectangle 12" shape, inside "Group 11".
# Split the single run into "storage" / "::" / "api" the same way the
# sibling "storage::datastore" label is split, by editing the trailing
# substrings in place (keeps the "storage::" prefix run untouched and
# replaces "manager" with "api").
$grp = $s.Shapes.Item(2)
$managerShape = $grp.GroupItems.Item(1)
$managerRange = $managerShape.TextFrame.TextRange

$sepRange = $managerRange.Characters(8, 2)
$sepRange.Text = "::"

$nameRange = $managerRange.Characters(10, 7)
$nameRange.Text = "api"

# --- *Manager -> *Storage rectangle labels -------------------------------
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $t = $shp.TextFrame.TextRange.Text
        if ($t -eq "EvaluationsManager") {
            $shp.TextFrame.TextRange.Text = "EvaluationsStorage"
        } elseif ($t -eq "AccountsManager") {
            $shp.TextFrame.TextRange.Text = "AccountsStorage"
        } elseif ($t -eq "CoursesManager") {
            $shp.TextFrame.TextRange.Text = "CoursesStorage"
        }
    }
}

# --- Reposition the connector under the storage::api box -----------------
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Straight Arrow Connector 41") {
        $shp.Left = 82.11035
    }
}
